# Update "想去人数" (want-to-go count) figures for the refreshed data pull.
# Affects both the "展览" sheet and the aggregated "全部类型" sheet, which
# mirrors the same events.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 3198
$wsExhibition.Range("F4").Value = 1019
$wsExhibition.Range("F5").Value = 304

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3198
$wsAll.Range("F4").Value = 1019
$wsAll.Range("F6").Value = 304
